# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# 1. Remove the extra "Student Ans"/"Correct Ans" blocks that lived in
#    columns G:H (rows 15-21) and in columns D:E for every row except
#    the header (15) and the two rows that now carry the real student
#    answers for those questions (16-17).
# -------------------------------------------------------------------
$ws.Range("G15:H21").Clear()
$ws.Range("D18:E40").Clear()

# -------------------------------------------------------------------
# 2. Re-style the row-10/11/12 labels (No., Marking, Total) to match
#    the "mtitleStyle" used by the header row above them (row 9).
# -------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A12").PasteSpecial(-4122)

# -------------------------------------------------------------------
# 3. Update the summary numbers (Right / Wrong / Not Attempt / Max,
#    Marking scheme, and Total score).
# -------------------------------------------------------------------
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 17
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 44
$ws.Range("E12").Value = "44/112"

# -------------------------------------------------------------------
# 4. Fill in the student answer column with the answers that match
#    the correct answer, using the same green "correctStyle" that the
#    correct-answer column already carries (copied from B10's style).
# -------------------------------------------------------------------
$ws.Range("B10").Copy()

$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Option A"

$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D17").Value = "Option C"

$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Option C"

$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Option D"

$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Option D"

$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = "Option D"

$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = "Option D"

$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "Option C"

$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("A35").Value = "Option D"

$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A38").Value = "Option A"

$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A40").Value = "Option D"
